# Fix mistake from unexposed prevalences in brm agg models
# Column E (Risk_Unexposed) previously held "NaN%" and should now hold
# the value that used to be in column F (Risk_Exposed). Column F should
# be updated with the newly computed exposed-risk percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @{ Row = 2; E = "14.9%"; F = "38.9%" },
    @{ Row = 3; E = "18.3%"; F = "41.5%" },
    @{ Row = 4; E = "16.8%"; F = "38.7%" },
    @{ Row = 5; E = "19.4%"; F = "40.7%" },
    @{ Row = 6; E = "21.2%"; F = "53.3%" },
    @{ Row = 7; E = "26.0%"; F = "54.3%" },
    @{ Row = 8; E = "23.8%"; F = "53.3%" },
    @{ Row = 9; E = "27.4%"; F = "53.2%" }
)

foreach ($item in $values) {
    $r = $item.Row
    $cellE = $ws.Cells.Item($r, 5)
    $cellF = $ws.Cells.Item($r, 6)
    $cellE.NumberFormat = "@"
    $cellF.NumberFormat = "@"
    $cellE.Value = $item.E
    $cellF.Value = $item.F
}

$wb.Save()
